$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new time-log entry row 26: date, hours, and task description (with wrap text),
# matching the formatting already used by the rows above it.
$ws.Range("A26").Value = 43519
$ws.Range("A26").NumberFormat = "d-mmm"

$ws.Range("B26").Value = 1.5

$ws.Range("D26").Value = "Indie Project: setting up test database, properties, copying Database class, creating sql for cleaning database before testing, starting UserDaoTest"
$ws.Range("D26").WrapText = $true
$ws.Rows.Item(26).RowHeight = 30

# Move the active selection to where the user left off editing.
$ws.Range("I29").Select() | Out-Null
